$wb = $excel.ActiveWorkbook

# Data for sheet "NBR" (sheet1): rows 3..17 -> A=1..15, B=6..20, C=values
$nbrC = @(856, 851, 816, 843, 845, 834, 828, 834, 826, 829, 817, 814, 814, 812, 822)

# Data for sheet "BAR" (sheet2): rows 3..17 -> A=1..15, B=6..20, C=values
$barC = @(921, 917, 954, 919, 906, 902, 900, 891, 888, 887, 896, 887, 887, 887, 861)

$wsNBR = $wb.Worksheets.Item("NBR")
$wsBAR = $wb.Worksheets.Item("BAR")

for ($i = 0; $i -lt 15; $i++) {
    $row = 3 + $i

    $wsNBR.Cells.Item($row, 1).Value = $i + 1
    $wsNBR.Cells.Item($row, 2).Value = $i + 6
    $wsNBR.Cells.Item($row, 3).Value = $nbrC[$i]

    $wsBAR.Cells.Item($row, 1).Value = $i + 1
    $wsBAR.Cells.Item($row, 2).Value = $i + 6
    $wsBAR.Cells.Item($row, 3).Value = $barC[$i]
}

# Apply the same formatting as A2 (border/bold/center) to the new A column cells
$wsNBR.Cells.Item(2, 1).Copy()
$wsNBR.Range($wsNBR.Cells.Item(3, 1), $wsNBR.Cells.Item(17, 1)).PasteSpecial(-4122)

$wsBAR.Cells.Item(2, 1).Copy()
$wsBAR.Range($wsBAR.Cells.Item(3, 1), $wsBAR.Cells.Item(17, 1)).PasteSpecial(-4122)
